$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$src = $ws.Range("A255:AF255")
$dst = $ws.Range("A300:AF300")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats = -4122
